$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the used range) used only to stage text values that
# look numeric (e.g. "305.53"), so pasting them into the data cells keeps them
# as text -- matching how those cells were already stored before this edit --
# instead of Excel auto-converting a plain numeric-looking assignment into a
# real number. Using a text FORMULA result + paste-special-values round trip
# avoids touching NumberFormat/styles, so it leaves no footprint behind.
$scratch = $ws.Range("ZZ1")


# Row 2
$ws.Range("D2").Value2 = "41.835.92"
$ws.Range("E2").Value2 = "  +4.28%  "

# Row 3
$ws.Range("D3").Value2 = "2.269.47"
$ws.Range("E3").Value2 = "  +2.24%  "

# Row 4
$ws.Range("E4").Value2 = "  -0.01%  "

# Row 5
$scratch.Formula = "=`"305.53`""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value2 = "  +4.00%  "

# Row 6
$scratch.Formula = "=`"92.21`""
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value2 = "  +5.31%  "

# Row 7
$ws.Range("E7").Value2 = "  +4.06%  "

# Row 8
$ws.Range("E8").Value2 = "  -0.04%  "

# Row 9
$ws.Range("E9").Value2 = "  +3.25%  "

# Row 10
$scratch.Formula = "=`"32.83`""
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value2 = "  +7.12%  "

# Row 11
$scratch.Formula = "=`"53.88`""
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value2 = "  +5.91%  "

# Row 12
$scratch.Formula = "=`"0.0798`""
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value2 = "  +2.25%  "

# Row 13
$ws.Range("E13").Value2 = "  +1.21%  "

# Row 14
$scratch.Formula = "=`"6.63`""
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value2 = "  +3.45%  "

# Row 15
$ws.Range("D15").Value2 = "2.623.27"
$ws.Range("E15").Value2 = "  +2.24%  "

# Row 16
$scratch.Formula = "=`"14.28`""
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value2 = "  +3.48%  "

# Row 17
$ws.Range("D17").Value2 = "2.293.75"
$ws.Range("E17").Value2 = "  +2.45%  "

# Row 18
$scratch.Formula = "=`"0.765`""
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value2 = "  +4.06%  "

# Row 19
$ws.Range("D19").Value2 = "41.771.06"
$ws.Range("E19").Value2 = "  +4.27%  "

# Row 20
$scratch.Formula = "=`"12.22`""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value2 = "  +8.60%  "

# Row 21
$ws.Range("D21").Value2 = "0.0₃0908"
$ws.Range("E21").Value2 = "  +2.12%  "

# Row 22
$ws.Range("E22").Value2 = "  +2.86%  "

# Row 23
$scratch.Formula = "=`"67.00`""
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value2 = "  +2.13%  "

# Row 24
$scratch.Formula = "=`"242.71`""
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value2 = "  +2.87%  "

# Row 25
$ws.Range("E25").Value2 = "  +5.02%  "

# Row 26
$ws.Range("E26").Value2 = "  +0.07%  "

# Row 27
$scratch.Formula = "=`"1.93`""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value2 = "  +5.80%  "

# Row 28
$scratch.Formula = "=`"24.19`""
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value2 = "  +4.59%  "

# Row 29
$ws.Range("B29").Value2 = "Toncoin"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$scratch.Formula = "=`"2.30`""
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value2 = "  +11.65%  "

# Row 30
$ws.Range("B30").Value2 = "Cosmos"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$scratch.Formula = "=`"9.62`""
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value2 = "  +3.23%  "

# Row 31
$scratch.Formula = "=`"34.24`""
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value2 = "  +7.77%  "

# Row 32
$scratch.Formula = "=`"158.34`""
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value2 = "  -0.34%  "

# Row 33
$scratch.Formula = "=`"1.00`""
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value2 = "  +0.01%  "

# Row 34
$ws.Range("E34").Value2 = "  +4.41%  "

# Row 35
$scratch.Formula = "=`"0.0749`""
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value2 = "  +4.98%  "

# Row 36
$ws.Range("E36").Value2 = "  +0.29%  "

# Row 37
$scratch.Formula = "=`"17.14`""
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value2 = "  +9.68%  "

# Row 38
$ws.Range("E38").Value2 = "  +1.75%  "

# Row 39
$ws.Range("E39").Value2 = "  +2.86%  "

# Row 40
$ws.Range("E40").Value2 = "  +4.86%  "

# Row 41
$scratch.Formula = "=`"1.81`""
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value2 = "  +2.86%  "

# Row 42
$ws.Range("E42").Value2 = "  +4.52%  "

# Row 43
$ws.Range("D43").Value2 = "2.069.55"
$ws.Range("E43").Value2 = "  -0.45%  "

# Row 44
$scratch.Formula = "=`"19.36`""
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value2 = "  +0.57%  "

# Row 45
$ws.Range("E45").Value2 = "  +3.46%  "

# Row 46
$scratch.Formula = "=`"10.33`""
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value2 = "  +3.18%  "

# Row 47
$ws.Range("E47").Value2 = "  +5.55%  "

# Row 48
$ws.Range("E48").Value2 = "  +7.21%  "

# Row 49
$ws.Range("B49").Value2 = "Stacks"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$scratch.Formula = "=`"1.52`""
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value2 = "  +3.63%  "

# Row 50
$ws.Range("B50").Value2 = "TrustWalletToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$scratch.Formula = "=`"1.16`""
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value2 = "  +3.34%  "

# Row 51
$ws.Range("B51").Value2 = "BitcoinSV"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$scratch.Formula = "=`"72.80`""
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value2 = "  +7.27%  "

# Clean up the scratch cell/clipboard so no trace of it is left behind
$scratch.Clear()
$excel.CutCopyMode = $false
